$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New naive QoQ errors observed for the latest ifoCAST vintage (one per row,
# appended as each row's matched-error series grows by one more quarter).
$newValues = @{
    2  = 0.3220726034063205
    3  = -1.732849868393072
    4  = 1.920156118130258
    5  = 1.480670296544459
    6  = -0.2710286406117954
    7  = 0.05367389642184558
    8  = 0.3403798556124878
    9  = 0.1895217986112106
    10 = -0.3589771035472806
}

$cols = @("B","C","D","E","F","G","H","I","J","K")

function Shift-Left($arr) {
    if ($arr.Count -le 1) {
        return @()
    }
    return $arr[1..($arr.Count - 1)]
}

for ($r = 2; $r -le 20; $r++) {
    # Read the current row's series (columns B..K), stopping at the first blank.
    $vals = @()
    foreach ($c in $cols) {
        $v2 = $ws.Range("$c$r").Value2
        if ($v2 -eq $null) {
            break
        }
        $vals += [double]$v2
    }

    # Drop the oldest (first) observation so everything shifts one column left.
    $shifted = Shift-Left $vals

    # Append the newly matched observation, if any, for this row.
    if ($newValues.ContainsKey($r)) {
        $shifted += [double]$newValues[$r]
    }

    # Write back the shifted series, then clear any now-unused trailing cells.
    for ($i = 0; $i -lt $cols.Count; $i++) {
        $c = $cols[$i]
        if ($i -lt $shifted.Count) {
            $ws.Range("$c$r").Value = $shifted[$i]
        } else {
            $ws.Range("$c$r").ClearContents()
        }
    }
}
